$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Address, $Value) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

$updates = @(
    @('D2', '44.142.98'),
    @('E2', '  +2.06%  '),
    @('D3', '2.389.28'),
    @('E3', '  +1.35%  '),
    @('E4', '  +0.12%  '),
    @('D5', '0.694'),
    @('E5', '  +7.04%  '),
    @('D6', '243.07'),
    @('E6', '  +4.19%  '),
    @('D7', '77.00'),
    @('E7', '  +7.66%  '),
    @('E8', '  +0.13%  '),
    @('D9', '0.639'),
    @('E9', '  +32.92%  '),
    @('D10', '0.103'),
    @('E10', '  +5.62%  '),
    @('D11', '57.69'),
    @('E11', '  +1.61%  '),
    @('D12', '33.65'),
    @('E12', '  +24.07%  '),
    @('D13', '7.59'),
    @('E13', '  +20.80%  '),
    @('E14', '  +2.13%  '),
    @('D15', '2.747.48'),
    @('E15', '  +1.26%  '),
    @('D16', '17.09'),
    @('E16', '  +5.82%  '),
    @('D17', '0.934'),
    @('E17', '  +7.95%  '),
    @('D18', '2.402.25'),
    @('E18', '  +2.11%  '),
    @('D19', '44.194.60'),
    @('E19', '  +2.08%  '),
    @('E20', '  +2.20%  '),
    @('D21', '6.71'),
    @('E21', '  +5.85%  '),
    @('D22', '78.61'),
    @('E22', '  +5.73%  '),
    @('D23', '258.93'),
    @('E23', '  +3.61%  '),
    @('E24', '  +0.04%  '),
    @('D25', '2.56'),
    @('E25', '  +4.34%  '),
    @('B26', 'Cosmos'),
    @('C26', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @('D26', '11.22'),
    @('E26', '  +12.12%  '),
    @('B27', 'WEMIXToken'),
    @('C27', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('D27', '3.68'),
    @('E27', '  -2.54%  '),
    @('D28', '1.77'),
    @('E28', '  +17.60%  '),
    @('B29', 'Toncoin'),
    @('C29', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('D29', '2.31'),
    @('E29', '  +5.77%  '),
    @('B30', 'EthereumClassic'),
    @('C30', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('D30', '23.46'),
    @('E30', '  +4.62%  '),
    @('D31', '175.88'),
    @('E31', '  +1.86%  '),
    @('D32', '0.131'),
    @('E32', '  +0.98%  '),
    @('D33', '0.136'),
    @('E33', '  +7.31%  '),
    @('D34', '5.39'),
    @('E34', '  +8.21%  '),
    @('D35', '0.0759'),
    @('E35', '  +10.15%  '),
    @('D36', '5.39'),
    @('E36', '  +6.50%  '),
    @('D37', '3.89'),
    @('E37', '  +5.42%  '),
    @('D38', '2.51'),
    @('E38', '  +3.24%  '),
    @('D39', '6.58'),
    @('E39', '  +0.41%  '),
    @('D40', '0.0277'),
    @('E40', '  +9.23%  '),
    @('B41', 'FraxShare'),
    @('C41', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('D41', '9.10'),
    @('E41', '  +2.06%  '),
    @('B42', 'InjectiveProtocol'),
    @('C42', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D42', '19.05'),
    @('E42', '  +2.51%  '),
    @('B43', 'Algorand'),
    @('C43', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'),
    @('D43', '0.202'),
    @('E43', '  +19.60%  '),
    @('E44', '  +0.01%  '),
    @('B45', 'NEARProtocol'),
    @('C45', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @('D45', '2.57'),
    @('E45', '  +15.91%  '),
    @('B46', 'ARBITRUM'),
    @('C46', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'),
    @('D46', '1.22'),
    @('E46', '  +4.92%  '),
    @('B47', 'Cronos'),
    @('C47', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D47', '0.102'),
    @('E47', '  +6.99%  '),
    @('B48', 'Aave'),
    @('C48', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D48', '103.98'),
    @('E48', '  +4.90%  '),
    @('B49', 'TrustWalletToken'),
    @('C49', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'),
    @('D49', '1.27'),
    @('E49', '  +5.39%  '),
    @('D50', '4.55'),
    @('E50', '  +1.95%  '),
    @('D51', '55.12'),
    @('E51', '  +9.89%  ')
)

foreach ($u in $updates) {
    Set-TextCell $ws $u[0] $u[1]
}
